# Daily attendance processing - 2025-12-12 17:27:48
# Reorders the "Recorded By" (column G) comma-separated list of recorders
# for specific rows so that the last-listed recorder is moved to the front.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$recordedByUpdates = @{
    2 = "System, backup@backdoor.com, system"
    3 = "System, dnasr281@gmail.com"
    4 = "System, backup@backdoor.com"
    5 = "System, backup@backdoor.com"
    6 = "System, dnasr281@gmail.com"
    8 = "System, backup@backdoor.com"
    10 = "System, dnasr281@gmail.com"
    11 = "System, dnasr281@gmail.com"
    12 = "System, dnasr281@gmail.com"
    13 = "System, dnasr281@gmail.com"
    14 = "System, dnasr281@gmail.com"
    15 = "System, dnasr281@gmail.com"
    17 = "System, dnasr281@gmail.com"
    18 = "System, dnasr281@gmail.com"
    19 = "System, dnasr281@gmail.com"
    20 = "System, dnasr281@gmail.com"
    21 = "System, dnasr281@gmail.com"
    22 = "System, dnasr281@gmail.com"
    24 = "System, dnasr281@gmail.com"
    26 = "System, dnasr281@gmail.com"
    28 = "System, backup@backdoor.com, system"
    29 = "System, dnasr281@gmail.com"
    30 = "System, backup@backdoor.com"
    31 = "System, backup@backdoor.com"
    32 = "System, dnasr281@gmail.com"
    34 = "System, backup@backdoor.com"
    36 = "System, dnasr281@gmail.com"
    37 = "System, dnasr281@gmail.com"
    38 = "System, dnasr281@gmail.com"
    39 = "System, dnasr281@gmail.com"
    40 = "System, dnasr281@gmail.com"
    41 = "System, dnasr281@gmail.com"
    43 = "System, dnasr281@gmail.com"
    44 = "System, dnasr281@gmail.com"
    45 = "System, dnasr281@gmail.com"
    46 = "System, dnasr281@gmail.com"
    47 = "System, dnasr281@gmail.com"
    48 = "System, dnasr281@gmail.com"
    50 = "System, dnasr281@gmail.com"
    52 = "System, dnasr281@gmail.com"
    54 = "System, backup@backdoor.com, system"
    55 = "System, dnasr281@gmail.com"
    56 = "System, backup@backdoor.com"
    57 = "System, backup@backdoor.com"
    58 = "System, dnasr281@gmail.com"
    60 = "System, backup@backdoor.com"
    62 = "System, dnasr281@gmail.com"
    63 = "System, dnasr281@gmail.com"
    64 = "System, dnasr281@gmail.com"
    65 = "System, dnasr281@gmail.com"
    66 = "System, dnasr281@gmail.com"
    67 = "System, dnasr281@gmail.com"
    69 = "System, dnasr281@gmail.com"
    70 = "System, dnasr281@gmail.com"
    71 = "System, dnasr281@gmail.com"
    72 = "System, dnasr281@gmail.com"
    73 = "System, dnasr281@gmail.com"
    74 = "System, dnasr281@gmail.com"
    76 = "System, dnasr281@gmail.com"
    78 = "System, dnasr281@gmail.com"
    80 = "System, backup@backdoor.com"
    81 = "System, backup@backdoor.com"
    82 = "System, backup@backdoor.com"
    83 = "System, dnasr281@gmail.com"
    84 = "System, dnasr281@gmail.com"
    85 = "System, dnasr281@gmail.com"
    86 = "System, dnasr281@gmail.com"
    87 = "admin@admin.com, dnasr281@gmail.com"
    90 = "System, dnasr281@gmail.com"
    92 = "System, dnasr281@gmail.com"
    93 = "System, dnasr281@gmail.com"
    94 = "System, dnasr281@gmail.com"
    96 = "System, dnasr281@gmail.com"
    99 = "System, dnasr281@gmail.com"
    101 = "System, dnasr281@gmail.com"
    106 = "System, backup@backdoor.com"
    107 = "System, backup@backdoor.com"
    108 = "System, backup@backdoor.com"
    109 = "System, dnasr281@gmail.com"
    110 = "System, dnasr281@gmail.com"
    111 = "System, dnasr281@gmail.com"
    112 = "System, dnasr281@gmail.com"
    113 = "admin@admin.com, dnasr281@gmail.com"
    116 = "System, dnasr281@gmail.com"
    118 = "System, dnasr281@gmail.com"
    119 = "System, dnasr281@gmail.com"
    120 = "System, dnasr281@gmail.com"
    122 = "System, dnasr281@gmail.com"
    125 = "System, dnasr281@gmail.com"
    127 = "System, dnasr281@gmail.com"
    132 = "System, backup@backdoor.com"
    133 = "System, backup@backdoor.com"
    134 = "System, backup@backdoor.com"
    135 = "System, dnasr281@gmail.com"
    136 = "System, dnasr281@gmail.com"
    137 = "System, dnasr281@gmail.com"
    138 = "System, dnasr281@gmail.com"
    139 = "admin@admin.com, dnasr281@gmail.com"
    142 = "System, dnasr281@gmail.com"
    144 = "System, dnasr281@gmail.com"
    145 = "System, dnasr281@gmail.com"
    146 = "System, dnasr281@gmail.com"
    148 = "System, dnasr281@gmail.com"
    151 = "System, dnasr281@gmail.com"
    153 = "System, dnasr281@gmail.com"
}

foreach ($rowNum in $recordedByUpdates.Keys) {
    $ws.Cells.Item($rowNum, 7).Value = $recordedByUpdates[$rowNum]
}
